$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405 (shifts existing rows 405-463 down to 406-464)
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with the weekly price-report entry
$ws.Range("A405").Value = 5
$ws.Range("B405").Value = "Macroferia Regional de Talca"
$ws.Range("C405").Value = "Maule"
$ws.Range("D405").Value = 44951
$ws.Range("E405").Value = 7
$ws.Range("F405").Value = 100114014
$ws.Range("G405").Value = "Betarraga"
$ws.Range("H405").Value = "Sin especificar"
$ws.Range("I405").Value = "Primera"
$ws.Range("J405").Value = 5000
$ws.Range("K405").Value = 700
$ws.Range("L405").Value = 700
$ws.Range("M405").Value = 700
$ws.Range("N405").Value = "$/paquete 5 unidades"
$ws.Range("O405").Value = "Región del Maule"
$ws.Range("P405").Value = 140
$ws.Range("Q405").Value = 5
$ws.Range("R405").Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D
$ws.Range("D405").NumberFormat = "YYYY-MM-DD HH:MM:SS"
